$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.95
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("V4").Value = 1.62
$ws.Range("M5").Value = 1.11
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("M6").Value = 1.08
$ws.Range("O6").Value = 1.36
$ws.Range("O10").Value = 1.12
$ws.Range("S10").Value = 1.26
$ws.Range("T10").Value = 3.5
$ws.Range("Q15").Value = 1.53
$ws.Range("U15").Value = 1.5
$ws.Range("Q16").Value = 1.6
$ws.Range("U16").Value = 1.53
$ws.Range("V16").Value = 2.38
$ws.Range("U17").Value = 1.91
$ws.Range("V17").Value = 1.8
$ws.Range("Q18").Value = 1.44
$ws.Range("Q19").Value = 1.7
$ws.Range("Q20").Value = 1.7
$ws.Range("Q21").Value = 1.48
$ws.Range("Q22").Value = 1.33
$ws.Range("R23").Value = 1.62
$ws.Range("G25").Value = 1.8
$ws.Range("M25").Value = 1.05
$ws.Range("O25").Value = 1.29
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.95
$ws.Range("M26").Value = 1.04
$ws.Range("O26").Value = 1.22
$ws.Range("U26").Value = 1.62
$ws.Range("G27").Value = 1.85
$ws.Range("M27").Value = 1.07
$ws.Range("O27").Value = 1.4
$ws.Range("V27").Value = 1.73
$ws.Range("M28").Value = 1.05
$ws.Range("O28").Value = 1.29
$ws.Range("U28").Value = 1.73
$ws.Range("M30").Value = 1.06
$ws.Range("O30").Value = 1.29
$ws.Range("U30").Value = 1.8
$ws.Range("V30").Value = 1.91
$ws.Range("R34").Value = 1.65
$ws.Range("U34").Value = 1.91
$ws.Range("V34").Value = 1.8
$ws.Range("U35").Value = 1.73
$ws.Range("R37").Value = 1.57
$ws.Range("V37").Value = 1.73
